$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet "Metadata" (sheet1): insert a new "Branch" column before the existing
# "Closing Balance" column (this shifts old column F -> G, carrying its
# header style along), then fill in the new column and tweak the other
# changed values.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Insert a new column at F; this shifts the old "Closing Balance" column (and
# its header formatting) from F to G, matching the diff exactly.
$wsMeta.Columns("F:F").Insert()

# New column header + data (plain text, Excel does not reinterpret these)
$wsMeta.Range("F1").Value = "Branch"
$wsMeta.Range("F2").Value = "MG Road, Bangalore"

# Update existing values that changed (plain text; not numeric/date-like
# enough for Excel to auto-convert them)
$wsMeta.Range("D2").Value = "01/05/2024 to 31/05/2024"
$wsMeta.Range("E2").Value = "1,25,000.00"

# ---------------------------------------------------------------------------
# Sheet "Transactions" (sheet2): add a new "Balance" column (E) with values,
# normalize date formats, and strip the currency symbol from Amount.
# ---------------------------------------------------------------------------
$wsTx = $wb.Worksheets.Item("Transactions")

# New header (plain text)
$wsTx.Range("E1").Value = "Balance"

# Some replacement strings look like genuine dates/numbers to Excel and would
# otherwise be silently converted into a date serial / numeric value with an
# auto-picked number format. Force those particular cells to be entered as
# plain text by temporarily marking them as Text, then restore an unstyled
# look by pasting the formatting from an always-unstyled neighbor cell (B-
# column "Description" cells never carry a style in this sheet).
$wsTx.Range("A2").NumberFormat = "@"
$wsTx.Range("A2").Value = "02/05/2024"
$wsTx.Range("B2").Copy()
$wsTx.Range("A2").PasteSpecial($xlPasteFormats)

$wsTx.Range("C2").NumberFormat = "@"
$wsTx.Range("C2").Value = "15,000.00"
$wsTx.Range("B2").Copy()
$wsTx.Range("C2").PasteSpecial($xlPasteFormats)

$wsTx.Range("E2").Value = "1,10,000.00"

$wsTx.Range("A3").NumberFormat = "@"
$wsTx.Range("A3").Value = "03/05/2024"
$wsTx.Range("B3").Copy()
$wsTx.Range("A3").PasteSpecial($xlPasteFormats)

$wsTx.Range("E3").NumberFormat = "@"
$wsTx.Range("E3").Value = "45,000.00"
$wsTx.Range("B3").Copy()
$wsTx.Range("E3").PasteSpecial($xlPasteFormats)

$wsTx.Range("E4").Value = "2,50,000.00"

$wsTx.Range("A5").NumberFormat = "@"
$wsTx.Range("A5").Value = "10/05/2024"
$wsTx.Range("B5").Copy()
$wsTx.Range("A5").PasteSpecial($xlPasteFormats)

$wsTx.Range("E5").NumberFormat = "@"
$wsTx.Range("E5").Value = "73,000.00"
$wsTx.Range("B5").Copy()
$wsTx.Range("E5").PasteSpecial($xlPasteFormats)

# Give the new "Balance" header the same style as the rest of row 1
$wsTx.Range("D1").Copy()
$wsTx.Range("E1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0
